# Refresh crypto price/volume snapshot (Price column D, Volume(1h) column E,
# plus a rank swap between WrappedEther/ShibaInu at rows 17-18).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'59.013.81"
$ws.Range("E2").Value = "  +1.36%  "
$ws.Range("D3").Value = "'2.587.67"
$ws.Range("E3").Value = "  -0.38%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'527.36"
$ws.Range("E5").Value = "  +0.92%  "
$ws.Range("D6").Value = "'139.10"
$ws.Range("E6").Value = "  -3.11%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  -0.84%  "
$ws.Range("D9").Value = "'2.599.33"
$ws.Range("E9").Value = "  -0.68%  "
$ws.Range("E10").Value = "  -0.72%  "
$ws.Range("E11").Value = "  -0.09%  "
$ws.Range("E12").Value = "  -3.67%  "
$ws.Range("E13").Value = "  +2.92%  "
$ws.Range("D14").Value = "'3.045.02"
$ws.Range("E14").Value = "  -0.37%  "
$ws.Range("D15").Value = "'58.949.61"
$ws.Range("E15").Value = "  +1.31%  "
$ws.Range("D16").Value = "'20.51"
$ws.Range("E16").Value = "  +0.58%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "'2.592.12"
$ws.Range("E17").Value = "  +1.12%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "'0.0000133"
$ws.Range("E18").Value = "  -0.97%  "
$ws.Range("D19").Value = "'344.36"
$ws.Range("E19").Value = "  +1.13%  "
$ws.Range("E20").Value = "  -0.67%  "
$ws.Range("E21").Value = "  -1.75%  "
$ws.Range("E22").Value = "  -0.62%  "
$ws.Range("D23").Value = "'0.998"
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("D24").Value = "'66.57"
$ws.Range("E24").Value = "  +1.73%  "
$ws.Range("E25").Value = "  -0.19%  "
$ws.Range("D26").Value = "'0.405"
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("E27").Value = "  +0.27%  "
$ws.Range("E28").Value = "  +0.36%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("E30").Value = "  -3.30%  "
$ws.Range("E31").Value = "  +1.50%  "
$ws.Range("D32").Value = "'5.91"
$ws.Range("E32").Value = "  -3.79%  "
$ws.Range("E33").Value = "  -0.49%  "
$ws.Range("D34").Value = "'149.58"
$ws.Range("E34").Value = "  -0.09%  "
$ws.Range("E35").Value = "  -0.97%  "
$ws.Range("E36").Value = "  -1.10%  "
$ws.Range("D37").Value = "'36.82"
$ws.Range("E37").Value = "  +2.20%  "
$ws.Range("E38").Value = "  +1.63%  "
$ws.Range("D39").Value = "'0.828"
$ws.Range("E39").Value = "  -4.33%  "
$ws.Range("D40").Value = "'0.812"
$ws.Range("E40").Value = "  -6.46%  "
$ws.Range("E41").Value = "  -0.57%  "
$ws.Range("D42").Value = "'0.998"
$ws.Range("E42").Value = "  +0.15%  "
$ws.Range("E43").Value = "  -1.28%  "
$ws.Range("D44").Value = "'269.41"
$ws.Range("E44").Value = "  -0.13%  "
$ws.Range("E45").Value = "  +0.75%  "
$ws.Range("E46").Value = "  -0.34%  "
$ws.Range("D47").Value = "'0.0515"
$ws.Range("E47").Value = "  -1.47%  "
$ws.Range("E48").Value = "  -1.91%  "
$ws.Range("D49").Value = "'1.962.86"
$ws.Range("E49").Value = "  -0.38%  "
$ws.Range("E50").Value = "  -0.01%  "
$ws.Range("D51").Value = "'18.23"
$ws.Range("E51").Value = "  -2.88%  "
